# Append a new price row (row 24) to the driver_prices sheet, mirroring
# the previous row's layout/format but with a fresh timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24
$prevRow = $newRow - 1

# Timestamp (column A) - carries the same date/time number format as the
# row above it.
$ws.Cells.Item($newRow, 1).Value2 = 44042.91666666666
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

# Driver prices (columns B through V), column N (K. Raikkonen) intentionally
# left blank for this update, same as the prior row.
$ws.Cells.Item($newRow, 2).Value2  = 31.2
$ws.Cells.Item($newRow, 3).Value2  = 29.7
$ws.Cells.Item($newRow, 4).Value2  = 26.1
$ws.Cells.Item($newRow, 5).Value2  = 23.7
$ws.Cells.Item($newRow, 6).Value2  = 20.9
$ws.Cells.Item($newRow, 7).Value2  = 20.9
$ws.Cells.Item($newRow, 8).Value2  = 15.5
$ws.Cells.Item($newRow, 9).Value2  = 13.7
$ws.Cells.Item($newRow, 10).Value2 = 12.9
$ws.Cells.Item($newRow, 11).Value2 = 12.2
$ws.Cells.Item($newRow, 12).Value2 = 10.2
$ws.Cells.Item($newRow, 13).Value2 = 9.9
$ws.Cells.Item($newRow, 15).Value2 = 9.6
$ws.Cells.Item($newRow, 16).Value2 = 9.699999999999999
$ws.Cells.Item($newRow, 17).Value2 = 8.6
$ws.Cells.Item($newRow, 18).Value2 = 7.8
$ws.Cells.Item($newRow, 19).Value2 = 5.9
$ws.Cells.Item($newRow, 20).Value2 = 6.1
$ws.Cells.Item($newRow, 21).Value2 = 5.7
$ws.Cells.Item($newRow, 22).Value2 = 9.699999999999999
